# Fruta / hortaliza, semanal
# Insert two new weekly rows (Florida King - Primera / Segunda) above the
# existing "Flavor Crest" row 67, shifting rows 67:79 down to 69:81.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows at row 67 (pushes old rows 67-79 down to 69-81,
# inheriting the date-formatted style from row 66 for column D, same as
# the target file).
$ws.Range("A67:A68").EntireRow.Insert()

# New row 67: Comercializadora del Agro de Limarí - Florida King - Primera
$ws.Range("A67").Value = 2
$ws.Range("B67").Value = 'Comercializadora del Agro de Limarí'
$ws.Range("C67").Value = 'Coquimbo'
$ws.Range("D67").Value = 44510
$ws.Range("E67").Value = 4
$ws.Range("F67").Value = 'Fruta'
$ws.Range("G67").Value = 100103
$ws.Range("H67").Value = 'Frutos de hueso (carozo)'
$ws.Range("I67").Value = 100103004
$ws.Range("J67").Value = 'Durazno'
$ws.Range("K67").Value = 'Florida King'
$ws.Range("L67").Value = 'Primera'
$ws.Range("M67").Value = 360
$ws.Range("N67").Value = 12000
$ws.Range("O67").Value = 13000
$ws.Range("P67").Value = 12500
$ws.Range("Q67").Value = '$/bandeja 10 kilos granel'
$ws.Range("R67").Value = 'Provincia de Limarí'
$ws.Range("S67").Value = 1250
$ws.Range("T67").Value = 10

# New row 68: Comercializadora del Agro de Limarí - Florida King - Segunda
$ws.Range("A68").Value = 2
$ws.Range("B68").Value = 'Comercializadora del Agro de Limarí'
$ws.Range("C68").Value = 'Coquimbo'
$ws.Range("D68").Value = 44510
$ws.Range("E68").Value = 4
$ws.Range("F68").Value = 'Fruta'
$ws.Range("G68").Value = 100103
$ws.Range("H68").Value = 'Frutos de hueso (carozo)'
$ws.Range("I68").Value = 100103004
$ws.Range("J68").Value = 'Durazno'
$ws.Range("K68").Value = 'Florida King'
$ws.Range("L68").Value = 'Segunda'
$ws.Range("M68").Value = 260
$ws.Range("N68").Value = 9000
$ws.Range("O68").Value = 10000
$ws.Range("P68").Value = 9500
$ws.Range("Q68").Value = '$/bandeja 10 kilos granel'
$ws.Range("R68").Value = 'Provincia de Limarí'
$ws.Range("S68").Value = 950
$ws.Range("T68").Value = 10
